$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting for cells holding numeric-looking overs values (e.g. "4.0" -> "3.0")
# so the engine keeps them as text instead of auto-converting to numbers.
$textCells = @("C16","L16","B21","B22","B23","B24","B25","K25")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("B2").Value = 9
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 'LBW'
$ws.Range("K2").Value = 27
$ws.Range("L2").Value = 12
$ws.Range("M2").Value = 'LBW'
$ws.Range("N2").Value = ' Hasan Ali'

# Row 3
$ws.Range("B3").Value = 11
$ws.Range("C3").Value = 5
$ws.Range("E3").Value = ' Dushmantha Chameera'
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 2
$ws.Range("N3").Value = ' Hasan Ali'

# Row 4
$ws.Range("B4").Value = 13
$ws.Range("C4").Value = 5
$ws.Range("D4").Value = 'Caught'
$ws.Range("E4").Value = ' Maheesh Theekshana'
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 'Caught'

# Row 5
$ws.Range("B5").Value = 20
$ws.Range("C5").Value = 9
$ws.Range("D5").Value = 'Caught'
$ws.Range("E5").Value = ' Nuwan Pradeep'
$ws.Range("K5").Value = 16
$ws.Range("L5").Value = 7
$ws.Range("M5").Value = 'NOT OUT'
$ws.Range("N5").Value = ' '

# Row 6
$ws.Range("B6").Value = 33
$ws.Range("C6").Value = 12
$ws.Range("D6").Value = 'Bowled'
$ws.Range("E6").Value = ' Nuwan Pradeep'
$ws.Range("K6").Value = 13
$ws.Range("L6").Value = 6
$ws.Range("M6").Value = 'LBW'

# Row 7
$ws.Range("B7").Value = 7
$ws.Range("C7").Value = 3
$ws.Range("D7").Value = 'Caught'
$ws.Range("E7").Value = ' Maheesh Theekshana'
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 'Bowled'
$ws.Range("N7").Value = ' Imad Wasim'

# Row 8
$ws.Range("B8").Value = 45
$ws.Range("C8").Value = 16
$ws.Range("D8").Value = 'LBW'
$ws.Range("E8").Value = ' Nuwan Pradeep'
$ws.Range("K8").Value = 5
$ws.Range("L8").Value = 4
$ws.Range("M8").Value = 'Bowled'
$ws.Range("N8").Value = ' Shadab Khan'

# Row 9
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 'Bowled'
$ws.Range("E9").Value = ' Dushmantha Chameera'
$ws.Range("K9").Value = 26
$ws.Range("L9").Value = 9
$ws.Range("M9").Value = 'Bowled'
$ws.Range("N9").Value = ' Shaheen Afridi'

# Row 10
$ws.Range("B10").Value = 10
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 'Bowled'
$ws.Range("E10").Value = ' Chamika Karunarathne'
$ws.Range("K10").Value = 29
$ws.Range("L10").Value = 9
$ws.Range("M10").Value = 'Caught'

# Row 11
$ws.Range("B11").Value = 57
$ws.Range("C11").Value = 17
$ws.Range("D11").Value = 'LBW'
$ws.Range("E11").Value = ' Chamika Karunarathne'
$ws.Range("K11").Value = 12
$ws.Range("L11").Value = 3
$ws.Range("M11").Value = 'Bowled'
$ws.Range("N11").Value = ' Imad Wasim'

# Row 12
$ws.Range("B12").Value = 24
$ws.Range("C12").Value = 8
$ws.Range("D12").Value = 'NOT OUT'
$ws.Range("K12").Value = 5
$ws.Range("L12").Value = 3
$ws.Range("M12").Value = 'LBW'

# Row 16
$ws.Range("A16").Value = 229
$ws.Range("B16").Value = 10
$ws.Range("C16").Value = '14.0'
$ws.Range("D16").Value = 84
$ws.Range("J16").Value = 136
$ws.Range("L16").Value = '9.3'
$ws.Range("M16").Value = 57

# Row 21
$ws.Range("A21").Value = 'Wanindu Hasaranga'
$ws.Range("B21").Value = '2.0'
$ws.Range("C21").Value = 37
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 18.5
$ws.Range("M21").Value = 0

# Row 22
$ws.Range("A22").Value = 'Nuwan Pradeep'
$ws.Range("B22").Value = '3.0'
$ws.Range("C22").Value = 36
$ws.Range("D22").Value = 4
$ws.Range("E22").Value = 12
$ws.Range("L22").Value = 28
$ws.Range("M22").Value = 1
$ws.Range("N22").Value = 14

# Row 23
$ws.Range("A23").Value = 'Maheesh Theekshana'
$ws.Range("B23").Value = '3.0'
$ws.Range("C23").Value = 58
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 19.33
$ws.Range("L23").Value = 30
$ws.Range("N23").Value = 15

# Row 24
$ws.Range("A24").Value = 'Dushmantha Chameera'
$ws.Range("B24").Value = '3.0'
$ws.Range("C24").Value = 50
$ws.Range("D24").Value = 2
$ws.Range("E24").Value = 16.67
$ws.Range("L24").Value = 29
$ws.Range("M24").Value = 4
$ws.Range("N24").Value = 14.5

# Row 25
$ws.Range("A25").Value = 'Chamika Karunarathne'
$ws.Range("B25").Value = '3.0'
$ws.Range("C25").Value = 48
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 16
$ws.Range("K25").Value = '1.3'
$ws.Range("L25").Value = 20
$ws.Range("N25").Value = 15.38
